$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is unambiguous (already won't be auto-parsed as a
# number by Excel: percentages, names, URLs, and "price" strings that
# contain more than one '.' e.g. "22.484.51"). These can be written with a
# plain .Value assignment.
$plainUpdates = @{
    "D2"  = "22.484.51"; "E2"  = "  +0.55%  ";
    "D3"  = "1.573.60";  "E3"  = "  +0.44%  ";
    "E4"  = "  -0.08%  ";
    "E5"  = "  -0.03%  ";
    "E6"  = "  +0.20%  ";
    "E7"  = "  -2.04%  ";
    "E8"  = "  +1.38%  ";
    "E9"  = "  -0.47%  ";
    "E10" = "  +0.38%  ";
    "E11" = "  -0.74%  ";
    "E12" = "  -0.07%  ";
    "E13" = "  +1.05%  ";
    "E14" = "  +0.66%  ";
    "E15" = "  +0.59%  ";
    "D16" = "1.576.62"; "E16" = "  +0.55%  ";
    "E17" = "  -1.16%  ";
    "E18" = "  +0.86%  ";
    "E19" = "  +0.28%  ";
    "E20" = "  -0.05%  ";
    "E21" = "  +1.76%  ";
    "E22" = "  -1.00%  ";
    "E23" = "  +2.34%  ";
    "D24" = "22.494.37"; "E24" = "  +0.57%  ";
    "E25" = "  -1.56%  ";
    "E26" = "  -3.12%  ";
    "E27" = "  -0.09%  ";
    "E28" = "  +1.07%  ";
    "E29" = "  +0.37%  ";
    "D31" = "1.750.84"; "E31" = "  +0.57%  ";
    "E32" = "  +7.92%  ";
    "E33" = "  +2.71%  ";
    "E34" = "  -0.12%  ";
    "E35" = "  -3.39%  ";
    "E36" = "  -1.24%  ";
    "E37" = "  -0.72%  ";
    "E38" = "  -4.17%  ";
    "E39" = "  +0.33%  ";
    "E40" = "  +0.91%  ";
    "E41" = "  +0.78%  ";
    "E42" = "  +0.24%  ";
    "E43" = "  -1.45%  ";
    "B44" = "Frax"; "C44" = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"; "E44" = "  -0.04%  ";
    "B45" = "EnergySwap"; "C45" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; "E45" = "  +0.67%  ";
    "E46" = "  +0.04%  ";
    "E48" = "  +3.66%  ";
    "E49" = "  -0.49%  ";
    "E50" = "  -2.72%  ";
    "E51" = "  +0.03%  ";
}

# Cells whose new text would otherwise be parsed as a plain number by Excel
# (single decimal point, e.g. "291.53"). Force them to stay text, matching
# the source data which stores every price as a string.
$textUpdates = @{
    "D6"  = "291.53";
    "D7"  = "0.3704";
    "D8"  = "49.90";
    "D9"  = "0.3386";
    "D10" = "1.143";
    "D11" = "0.07547";
    "D13" = "21.28";
    "D14" = "6.025";
    "D15" = "6.961";
    "D18" = "90.70";
    "D19" = "0.06759";
    "D20" = "1.001";
    "D21" = "6.304";
    "D22" = "16.44";
    "D25" = "2.362";
    "D26" = "2.606";
    "D27" = "20.09";
    "D28" = "149.11";
    "D29" = "5.049";
    "D30" = "125.24";
    "D33" = "6.249";
    "D34" = "2.013";
    "D35" = "9.773";
    "D36" = "0.08354";
    "D37" = "0.02491";
    "D38" = "1.363";
    "D39" = "0.2302";
    "D40" = "0.06547";
    "D41" = "5.448";
    "D42" = "11.36";
    "D43" = "0.6230";
    "D44" = "1.001";
    "D45" = "14.06";
    "D47" = "0.5858";
    "D48" = "129.28";
    "D50" = "1.223";
    "D51" = "0.07326";
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

foreach ($ref in $textUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$ref]
}
